# Cafe_Data_20_rows_excel_check.xlsx
# Commit: "Added API route for profit margin (KPI)"
#
# 1. Tidy up the "Total revenue" label (drop the trailing colon).
# 2. Add a new KPI row (30): Profit margin (%) = (Total Profit / Total Revenue) * 100,
#    with a value formula and a "same result" boolean-check formula, mirroring the
#    three rows above it (Total revenue / Total cost / Total profit).
# 3. Nudge the selection / column B width to match the post-edit sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the "Total revenue" label text (remove trailing colon) ---
$ws.Range("A27").Value = "Total revenue = sum (rate * quantity + tax " + [char]8211 + " discount)"

# --- 2. Add the new Profit margin (%) KPI row ---
$ws.Range("A30").Value = "Profit margin (%) = (Total Profit / Total Revenue) * 100"
$ws.Range("B30").Formula = "=(B29/B27) * 100"
$ws.Range("C30").Formula = "=32.7663801037014=B30"

# Match the formatting used by the row above it (Total profit row)
$ws.Range("A29:C29").Copy() | Out-Null
$ws.Range("A30:C30").PasteSpecial(-4122) | Out-Null

# --- 3. Misc view/layout tweaks to mirror the post-edit workbook ---
$ws.Range("C31").Select() | Out-Null
$ws.Columns.Item(2).ColumnWidth = 10.3
